# EI Variable Installments T2 scenarios
# Insert a new "waittopageload1" / 2000 row above the "Submit" row on the
# "Edit Repayment Schedule" sheet, and make that sheet the active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edit Repayment Schedule")

# Insert a new row above row 11 ("clickonsubmit" / "Submit"), shifting the
# remaining rows down by one.
$ws.Rows.Item(11).Insert()

# Populate the new row with the waittopageload1 step and its value.
$ws.Cells.Item(11, 1).Value = "waittopageload1"
$ws.Cells.Item(11, 2).Value = 2000

# Match the formatting used by the other "Pattern" value rows (e.g. B3),
# which use a different cell style than the row that used to occupy row 11.
$ws.Cells.Item(3, 2).Copy()
$ws.Cells.Item(11, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Make "Edit Repayment Schedule" the active sheet/tab, and select A11:B11
# on it, which also clears the previous tab-selected state on NewLoanInput.
$ws.Activate()
$ws.Range("A11:B11").Select()
